$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "100 chart game",
    "12 string neck",
    "123 stores",
    "16w capris",
    "1st endurance",
    "2 inch hip bands",
    "20 30 compression leggings",
    "20 ways to draw a dress",
    "2018 club dresses for women",
    "2018 fitness tracker",
    "2x capri",
    "2x capri leggings",
    "2x capri leggings plus size",
    "2x compression leggings women",
    "2x womens compression shorts",
    "360 compression knee brace",
    "360 compression knee brace for men",
    "360 knee compression brace",
    "360 leggings",
    "360 tank top",
    "3x black stockings",
    "3x compression leggings",
    "3x compression pantyhose",
    "3x control top tights",
    "3x knee brace",
    "3x knee support",
    "3x stockings",
    "3x womens leggings",
    "3x3 fitness rings",
    "4 black",
    "4 target",
    "4 ways stretch leggings",
    "4 women",
    "4 x",
    "4x back brace support",
    "4x black leggings",
    "4x compression shirt",
    "4x compression shirts",
    "4x knee braces",
    "4x knee support",
    "4x neoprene pants",
    "4x workout pants for women",
    "4x4 recovery equipment",
    "4x4 recovery gear",
    "5 xl knee brace",
    "5x black tights",
    "5x women clothing",
    "6 pack of leggings for women",
    "7 compression shorts women",
    "7 day flat belly",
    "7 inch basketball balls",
    "7 inch compression shorts",
    "7 inch nails",
    "7 inch running shorts women",
    "7 joint support",
    "7 leggings for women",
    "7 pack leggings",
    "7 pack leggings for women",
    "7 pack of leggings",
    "7 pair compression",
    "7 pair leggings",
    "7 pair of leggings",
    "7 pairs of leggings",
    "70 30 waist trainer",
    "80 awesome sports games",
    "80 s leggings",
    "80 styles clothing for women",
    "80 workout clothes for women",
    "90 clothing",
    "90 degree capri",
    "90 degrees capri",
    "a black mile to the surface",
    "a game of nasty things",
    "a nasty bit of rough",
    "a time to rise",
    "a400 men",
    "ab bike",
    "ab core and pelvic",
    "ab cuts",
    "ab freezer",
    "ab generator",
    "ab pro",
    "ab pro 360",
    "abdomen compression support",
    "abdomen compression women",
    "abdomen support",
    "abdominal compression for women",
    "abdominal compression garment",
    "abdominal cuts",
    "abdominal machine",
    "abdominal support band",
    "abdominal supports",
    "abdominal training",
    "abisko lite trekking",
    "abisko trekking",
    "about time nighttime recovery",
    "abs for women",
    "abs legging",
    "absolute rings for women",
    "accelerate compression shorts"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
